$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: move a row's data (values/formulas) + cell formatting from one row
# number to another, leaving the source row's formatting/values behind (it
# gets overwritten later by whatever ends up being placed on top of it).
# ---------------------------------------------------------------------------
function Move-RowData($ws, $srcRow, $dstRow, $lastCol) {
    $srcRange = $ws.Range("A" + $srcRow + ":" + $lastCol + $srcRow)
    $dstRange = $ws.Range("A" + $dstRow + ":" + $lastCol + $dstRow)
    $n = $srcRange.Columns.Count
    for ($i = 1; $i -le $n; $i++) {
        $srcCell = $srcRange.Cells.Item(1, $i)
        $dstCell = $dstRange.Cells.Item(1, $i)
        $dstCell.Formula = $srcCell.Formula
    }
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)   # xlPasteFormats
    $ws.Application.CutCopyMode = 0
}

# ---------------------------------------------------------------------------
# Relocate the tail of the journal (rows 50-56) six rows down to 56-62 so we
# can open up space for the six new journal entries. Destinations are filled
# from the bottom up so that nothing is overwritten before it is copied.
# ---------------------------------------------------------------------------
Move-RowData $ws 56 62 "C"

# Row 62's total formula now spans the new range; fix it up (and restore its
# pristine "General" number format) right away, before row 56 gets reused.
$ws.Range("C62").Formula = "=SUM(C4:C61)"
$ws.Range("C56").Copy()
$ws.Range("C62").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

Move-RowData $ws 55 61 "D"
Move-RowData $ws 54 60 "D"
Move-RowData $ws 53 59 "D"
Move-RowData $ws 52 58 "D"
Move-RowData $ws 51 56 "D"
Move-RowData $ws 50 53 "D"

# thick-bottom rows keep their taller row height
$ws.Rows("61:61").RowHeight = 15.75
$ws.Rows("62:62").RowHeight = 15.75

# ---------------------------------------------------------------------------
# Stamp the standard data-row look (borders/wrap) used throughout the journal
# onto the freshly opened rows 50-52, 54-55 and 57, using row 49 (A/B/C/D =
# styles 7/11/13/17) as the template, then clear any stale values that were
# left behind by the moves above (column A in particular).
# ---------------------------------------------------------------------------
foreach ($r in 50, 51, 52, 54, 55, 57) {
    $ws.Range("A49:D49").Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = 0
    $ws.Range("A" + $r + ":D" + $r).ClearContents()
}

# Rows 55 and 56 used to be the last two (thick-bottom-bordered) rows of the
# table; now that they hold ordinary entries, drop the stale taller
# thick-bottom row height that tagged along with the row numbers.
$ws.Rows("55:55").AutoFit()
$ws.Rows("56:56").AutoFit()

# ---------------------------------------------------------------------------
# New journal entries (values). Shared strings must be created in this exact
# order so they land on unique-string indices 48, 49, 50.
# ---------------------------------------------------------------------------
$ws.Range("B50").Value2 = "Création du script de prise de vue pour le raspberry"
$ws.Range("C50").Value2 = 1.5

$ws.Range("B51").Value2 = "Debug"
$ws.Range("C51").Value2 = 2

$ws.Range("B52").Value2 = "Documentation"
$ws.Range("C52").Value2 = 1.5

$ws.Range("B54").Value2 = "Debug"
$ws.Range("C54").Value2 = 4
$ws.Range("D54").Value2 = "Problème de compatibilité relatifs aux formats d'image"

$ws.Range("B55").Value2 = "Documentation"
$ws.Range("C55").Value2 = 2.5

$ws.Range("B57").Value2 = "Documentation"
$ws.Range("C57").Value2 = 4
$ws.Range("D57").Value2 = "Commentaire du code"

Write-Output "done"
